# Regenerate merged AHB files
# 1) Rename the "_old" / "_new" suffixed header labels (row 1) to the
#    version-specific "_FV2410" / "_FV2504" suffixes used by the newly
#    merged AHB export.
# 2) Turn the sheet's data range into an Excel Table ("Table1") so it can
#    be filtered/sorted like the other merged AHB sheets.
# 3) Freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fields = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $fields.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fields[$i] + "_FV2410"
}

$ws.Cells.Item(1, 11).Value = "diff"

for ($i = 0; $i -lt $fields.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fields[$i] + "_FV2504"
}

# Convert the used range into a proper table.
$rng = $ws.UsedRange
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row.
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
